# Actualizar 05-12-2020 17-04-35
# Adds Choluteca-department health-facility rows to HOSPITALES_HN table and
# refreshes the Admin2 codes (J column) for the previously-added municipios
# in Copan (rows 169-173) to use their full 4-digit code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix up Admin2_id (column J) on rows 169-173: the 2-digit codes are
#    replaced with the full 4-digit municipio code (matches the pattern
#    used by every other row in the sheet, e.g. "0320", "0321", ...).
# ---------------------------------------------------------------------
$ws.Range("J169").Value = "0406"
$ws.Range("J170").Value = "0408"
$ws.Range("J171").Value = "0419"
$ws.Range("J172").Value = "0413"
$ws.Range("J173").Value = "0410"

# ---------------------------------------------------------------------
# 2) Append the new rows (174-185) to the bottom of the table, re-using
#    the formatting already present on the last existing data row (173)
#    so the new cells inherit the correct number formats / borders.
# ---------------------------------------------------------------------
$ws.Range("B173:S173").Copy()
$ws.Range("B174:S185").PasteSpecial(-4122)
$ws.Range("V173:W173").Copy()
$ws.Range("V174:V179").PasteSpecial(-4122)
$ws.Range("W173:W173").Copy()
$ws.Range("W174:W179").PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-Row($r, $vals) {
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}

Set-Row 174 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=12;
    J="0612"; K="San Antonio de Flores"; L="Municipio"; M=1; N="061201"; O="San Antonio de Flores"; P="Aldea";
    Q="HND-0612"; R="Salud"; S="Cesamo"; V=13.665829; W=-87.362111999999996
}
Set-Row 175 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=8;
    J="0608"; K="Morolica"; L="Municipio"; M=1; N="060801"; O="Morolica"; P="Aldea";
    Q="HND-0608"; R="Salud"; S="Cesamo"; V=13.568472999999999; W=-86.907604000000006
}
Set-Row 176 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=16;
    J="0616"; K="Santa Ana de Yusguare"; L="Municipio"; M=1; N="061601"; O="Santa Ana de Yusguare"; P="Aldea";
    Q="HND-0616"; R="Salud"; S="Cesamo"; V=13.293342000000001; W=-87.111866000000006
}
Set-Row 177 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=7;
    J="0607"; K="Marcovia"; L="Municipio"; M=1; N="060701"; O="Marcovia"; P="Aldea";
    Q="HND-0607"; R="Salud"; S="Cesamo"; V=13.284371999999999; W=-87.312591999999995
}
Set-Row 178 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=9;
    J="0609"; K="Namasigue"; L="Municipio"; M=1; N="060901"; O="Namasigue"; P="Aldea";
    Q="HND-0609"; R="Salud"; S="Cesamo"; V=13.203086000000001; W=-87.138760000000005
}
Set-Row 179 @{
    B="HND"; C="Honduras"; D=3; E=6; F="06"; G="Choluteca"; H="Departamento"; I=7;
    J="0607"; K="Marcovia"; L="Municipio"; M=14; N="060714"; O="Monjaras"; P="Aldea";
    Q="HND-0607"; R="Salud"; S="Cesamo"; V=13.198708; W=-87.375152999999997
}

# Rows 180-185: partially filled placeholder rows (category info only,
# awaiting the rest of the facility data) - E,F,J,N keep their formatted
# but empty cells from the paste above.
180..185 | ForEach-Object {
    $r = $_
    Set-Row $r @{
        B="HND"; C="Honduras"; D=3; H="Departamento"; L="Municipio"; P="Aldea";
        Q="HND-"; R="Salud"; S="Cesamo"
    }
}

# ---------------------------------------------------------------------
# 3) Grow the Excel table (ListObject) + AutoFilter to the new extent.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:W185"))

# ---------------------------------------------------------------------
# 4) Update the hidden _FilterDatabase defined name to match.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=HOSPITALES!`$A`$1:`$W`$185"
    }
}

# ---------------------------------------------------------------------
# 5) Update the view: frozen-pane scroll position + active selection.
# ---------------------------------------------------------------------
$aw = $excel.ActiveWindow
$ws.Range("W180").Select()
$aw.ScrollRow = 165
$aw.ScrollColumn = 20
